$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $value) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = "Normal"
}

Set-TextValue 2 4 '35.953.52'
$ws.Cells.Item(2, 5).Value = '  -1.87%  '

Set-TextValue 3 4 '1.994.11'
$ws.Cells.Item(3, 5).Value = '  -2.98%  '

$ws.Cells.Item(4, 5).Value = '  +0.00%  '

Set-TextValue 5 4 '246.63'
$ws.Cells.Item(5, 5).Value = '  -0.31%  '

Set-TextValue 6 4 '0.641'
$ws.Cells.Item(6, 5).Value = '  -3.47%  '

Set-TextValue 7 4 '59.65'
$ws.Cells.Item(7, 5).Value = '  +8.84%  '

$ws.Cells.Item(8, 5).Value = '  +0.02%  '

Set-TextValue 9 4 '58.78'
$ws.Cells.Item(9, 5).Value = '  -3.20%  '

$ws.Cells.Item(10, 5).Value = '  -0.25%  '

$ws.Cells.Item(11, 5).Value = '  -1.73%  '

$ws.Cells.Item(12, 5).Value = '  -2.03%  '

$ws.Cells.Item(13, 5).Value = '  -2.45%  '

Set-TextValue 14 4 '14.82'
$ws.Cells.Item(14, 5).Value = '  -0.11%  '

Set-TextValue 15 4 '2.283.47'
$ws.Cells.Item(15, 5).Value = '  -3.13%  '

$ws.Cells.Item(16, 5).Value = '  -2.44%  '

Set-TextValue 17 4 '19.56'
$ws.Cells.Item(17, 5).Value = '  +13.09%  '

Set-TextValue 18 4 '1.997.24'
$ws.Cells.Item(18, 5).Value = '  -2.99%  '

Set-TextValue 19 4 '35.851.18'
$ws.Cells.Item(19, 5).Value = '  -1.84%  '

Set-TextValue 20 4 '71.96'
$ws.Cells.Item(20, 5).Value = '  -0.29%  '

Set-TextValue 21 4 '0.0₃0852'
$ws.Cells.Item(21, 5).Value = '  -1.21%  '

Set-TextValue 22 4 '5.24'
$ws.Cells.Item(22, 5).Value = '  -0.36%  '

Set-TextValue 23 4 '233.67'
$ws.Cells.Item(23, 5).Value = '  -1.98%  '

$ws.Cells.Item(24, 2).Value = 'PancakeSwap'
$ws.Cells.Item(24, 3).Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue 24 4 '2.64'
$ws.Cells.Item(24, 5).Value = '  +15.62%  '

$ws.Cells.Item(25, 2).Value = 'Dai'
$ws.Cells.Item(25, 3).Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue 25 4 '1.00'
$ws.Cells.Item(25, 5).Value = '  +0.11%  '

$ws.Cells.Item(26, 5).Value = '  -4.29%  '

Set-TextValue 27 4 '9.66'
$ws.Cells.Item(27, 5).Value = '  +4.41%  '

Set-TextValue 28 4 '165.39'
$ws.Cells.Item(28, 5).Value = '  -0.55%  '

Set-TextValue 29 4 '19.43'
$ws.Cells.Item(29, 5).Value = '  -3.50%  '

$ws.Cells.Item(30, 5).Value = '  -1.86%  '

$ws.Cells.Item(31, 5).Value = '  -2.78%  '

$ws.Cells.Item(32, 5).Value = '  -6.41%  '

Set-TextValue 33 4 '0.0983'
$ws.Cells.Item(33, 5).Value = '  +13.76%  '

Set-TextValue 34 4 '0.0607'
$ws.Cells.Item(34, 5).Value = '  +1.80%  '

Set-TextValue 35 4 '2.50'
$ws.Cells.Item(35, 5).Value = '  +10.45%  '

Set-TextValue 36 4 '4.43'
$ws.Cells.Item(36, 5).Value = '  -1.42%  '

$ws.Cells.Item(37, 5).Value = '  +0.02%  '

$ws.Cells.Item(38, 5).Value = '  -1.39%  '

Set-TextValue 39 4 '5.75'
$ws.Cells.Item(39, 5).Value = '  +13.87%  '

$ws.Cells.Item(40, 5).Value = '  -1.45%  '

$ws.Cells.Item(41, 2).Value = 'HuobiToken'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue 41 4 '2.86'
$ws.Cells.Item(41, 5).Value = '  -1.24%  '

$ws.Cells.Item(42, 2).Value = 'VeChain'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 42 4 '0.0214'
$ws.Cells.Item(42, 5).Value = '  -0.64%  '

Set-TextValue 43 4 '0.0935'
$ws.Cells.Item(43, 5).Value = '  +1.85%  '

$ws.Cells.Item(44, 5).Value = '  +0.00%  '

$ws.Cells.Item(45, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue 45 4 '16.60'
$ws.Cells.Item(45, 5).Value = '  +3.78%  '

$ws.Cells.Item(46, 2).Value = 'Aave'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue 46 4 '94.25'
$ws.Cells.Item(46, 5).Value = '  -0.53%  '

$ws.Cells.Item(47, 5).Value = '  +3.29%  '

Set-TextValue 48 4 '1.371.85'
$ws.Cells.Item(48, 5).Value = '  -3.18%  '

$ws.Cells.Item(49, 2).Value = 'MXToken'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue 49 4 '2.90'
$ws.Cells.Item(49, 5).Value = '  -0.75%  '

$ws.Cells.Item(50, 2).Value = 'RenderToken'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue 50 4 '2.35'
$ws.Cells.Item(50, 5).Value = '  +2.96%  '

Set-TextValue 51 4 '46.92'
$ws.Cells.Item(51, 5).Value = '  +2.10%  '

